# ---------------------------------------------------------------------------
# Reproduces:
#   1. Slide 6's table changes its table style (tableStyleId) from the
#      deck-local "Table_0" style to the built-in PowerPoint table style
#      {4D7A8931-E6E0-40ED-B684-240378CCEAEE}.
#   2. The presentation's theme ("theme1.xml", used by the slide master and
#      therefore by every slide) swaps its 12-colour scheme (and effectively
#      its identity) from the "Integral" palette to the default "Office"
#      palette that used to live only in the otherwise-unused notes-master
#      theme ("theme2.xml").
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style on slide 6 (the table is the 2nd shape on that slide).
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{4D7A8931-E6E0-40ED-B684-240378CCEAEE}")

# ---------------------------------------------------------------------------
# 2) Theme colour swap: recolour the active theme ("Integral") to the
#    "Office" palette. RGB values must be supplied as the packed 0x00BBGGRR
#    integers PowerPoint's COM automation model expects.
#
#    Office palette (target) -> packed BGR integer:
#      dk1      000000 -> 0
#      lt1      FFFFFF -> 16777215
#      dk2      44546A -> 6968388
#      lt2      E7E6E6 -> 15132391
#      accent1  5B9BD5 -> 13998939
#      accent2  ED7D31 -> 3243501
#      accent3  A5A5A5 -> 10855845
#      accent4  FFC000 -> 49407
#      accent5  4472C4 -> 12874308
#      accent6  70AD47 -> 4697456
#      hlink    0563C1 -> 12673797
#      folHlink 954F72 -> 7491477
# ---------------------------------------------------------------------------
$themeColors = $p.Slides.Item(1).ThemeColorScheme

$themeColors.Item(1).RGB  = 0
$themeColors.Item(2).RGB  = 16777215
$themeColors.Item(3).RGB  = 6968388
$themeColors.Item(4).RGB  = 15132391
$themeColors.Item(5).RGB  = 13998939
$themeColors.Item(6).RGB  = 3243501
$themeColors.Item(7).RGB  = 10855845
$themeColors.Item(8).RGB  = 49407
$themeColors.Item(9).RGB  = 12874308
$themeColors.Item(10).RGB = 4697456
$themeColors.Item(11).RGB = 12673797
$themeColors.Item(12).RGB = 7491477
